$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing "source" block (rows 26-27) down to rows 32-33 ---
# (same text, same bold/italic formatting, just relocated further down the sheet
#  to make room for the new "Number of employees / Assets / Turnover" table)
$ws.Range("A32").Value = "SME Performance Review EU"
$ws.Range("A32").Font.Bold = $true

$ws.Range("A33").Value = "SME Performance Review EU, ""SBA Fact sheet"", 2013.  Available at http://ec.europa.eu/enterprise/policies/sme/facts-figures-analysis/performance-review/index_en.htm"
$ws.Range("A33").Font.Italic = $true

# --- New table: "Number of employees / Assets / Turnover" thresholds (rows 23-27) ---
$ws.Range("B23").Value = "Number of employees"
$ws.Range("B23").Font.Bold = $true
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C23").Font.Bold = $true
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D23").Font.Bold = $true

$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"

$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"

$ws.Range("A26").Value = "Medium"
$ws.Range("A26").Font.Bold = $false
$ws.Range("B26").Value = "<250"

$ws.Range("A27").Value = "Large"
$ws.Range("A27").Font.Italic = $false
$ws.Range("B27").Value = ">249"
